$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 86, shifting rows 86:169 down to 88:171.
$ws.Rows("86:87").Insert()

# Row 86 (new): Comercializadora del Agro de Limari - Pepino ensalada - Primera, 28-Apr-2022
$ws.Range("A86").Value = 2
$ws.Range("B86").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C86").Value = 'Coquimbo'
$ws.Range("D86").Value = 44679
$ws.Range("E86").Value = 4
$ws.Range("F86").Value = 100112043
$ws.Range("G86").Value = 'Pepino ensalada'
$ws.Range("H86").Value = 'Sin especificar'
$ws.Range("I86").Value = 'Primera'
$ws.Range("J86").Value = 600
$ws.Range("K86").Value = 13000
$ws.Range("L86").Value = 14000
$ws.Range("M86").Value = 13500
$ws.Range("N86").Value = '$/caja 60 unidades'
$ws.Range("O86").Value = 'Provincia de Limarí'
$ws.Range("P86").Value = 225
$ws.Range("Q86").Value = 60
$ws.Range("R86").Value = 'Hortaliza'

# Row 87 (new): Comercializadora del Agro de Limari - Pepino ensalada - Segunda, 28-Apr-2022
$ws.Range("A87").Value = 2
$ws.Range("B87").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C87").Value = 'Coquimbo'
$ws.Range("D87").Value = 44679
$ws.Range("E87").Value = 4
$ws.Range("F87").Value = 100112043
$ws.Range("G87").Value = 'Pepino ensalada'
$ws.Range("H87").Value = 'Sin especificar'
$ws.Range("I87").Value = 'Segunda'
$ws.Range("J87").Value = 200
$ws.Range("K87").Value = 11000
$ws.Range("L87").Value = 12000
$ws.Range("M87").Value = 11500
$ws.Range("N87").Value = '$/caja 100 unidades'
$ws.Range("O87").Value = 'Provincia de Limarí'
$ws.Range("P87").Value = 115
$ws.Range("Q87").Value = 100
$ws.Range("R87").Value = 'Hortaliza'
